{"js": "// Fix the typo \"Lites\" -> \"Liste\" in the \"-Lites des musiques\" table cell\n// (description of the \"available\" command), turning it into\n// \"-Liste des musiques\".\nconst body = context.document.body;\n\nconst results = body.search(\"-Lites des musiques\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find text \"-Lites des musiques\" to fix.');\n}\n\nresults.items[0].insertText(\"-Liste des musiques\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix the typo \"Lites\" -> \"Liste\" in the \"-Lites des musiques\" table cell\n# (description of the \"available\" command), turning it into\n# \"-Liste des musiques\".\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"-Lites des musiques\", $false, $false, $false, $false, $false, $true, 1, $false, \"-Liste des musiques\", 2)\n\nif (-not $found) {\n    throw 'Could not find text \"-Lites des musiques\" to fix.'\n}\n"}
